$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'25.761.43"
$ws.Range("E2").Value = "  -4.01%  "

# Row 3
$ws.Range("D3").Value = "'1.816.79"
$ws.Range("E3").Value = "  -3.06%  "

# Row 4
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").Value = "'278.85"

# Row 6
$ws.Range("D6").Value = "'1.000"

# Row 7
$ws.Range("D7").Value = "'0.5095"
$ws.Range("E7").Value = "  -4.79%  "

# Row 8
$ws.Range("D8").Value = "'0.3539"
$ws.Range("E8").Value = "  -5.54%  "

# Row 9
$ws.Range("D9").Value = "'44.78"

# Row 10
$ws.Range("D10").Value = "'0.06665"
$ws.Range("E10").Value = "  -7.37%  "

# Row 11
$ws.Range("D11").Value = "'20.10"
$ws.Range("E11").Value = "  -7.05%  "

# Row 12
$ws.Range("D12").Value = "'0.8273"
$ws.Range("E12").Value = "  -7.12%  "

# Row 13
$ws.Range("D13").Value = "'0.07907"
$ws.Range("E13").Value = "  -3.50%  "

# Row 14
$ws.Range("D14").Value = "'1.793.01"
$ws.Range("E14").Value = "  -4.41%  "

# Row 15
$ws.Range("D15").Value = "'5.080"
$ws.Range("E15").Value = "  -4.39%  "

# Row 16
$ws.Range("D16").Value = "'87.78"
$ws.Range("E16").Value = "  -5.92%  "

# Row 17
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.41%  "

# Row 18
$ws.Range("D18").Value = "'14.09"
$ws.Range("E18").Value = "  -5.03%  "

# Row 19
$ws.Range("D19").Value = "'0.000008031"
$ws.Range("E19").Value = "  -5.90%  "

# Row 20
$ws.Range("E20").Value = "  -0.17%  "

# Row 21
$ws.Range("D21").Value = "'25.806.73"
$ws.Range("E21").Value = "  -3.96%  "

# Row 22
$ws.Range("D22").Value = "'4.752"
$ws.Range("E22").Value = "  -4.83%  "

# Row 23
$ws.Range("D23").Value = "'9.992"
$ws.Range("E23").Value = "  -5.81%  "

# Row 24
$ws.Range("D24").Value = "'6.124"
$ws.Range("E24").Value = "  -4.15%  "

# Row 25
$ws.Range("D25").Value = "'2.227"
$ws.Range("E25").Value = "  -2.64%  "

# Row 26
$ws.Range("D26").Value = "'142.34"

# Row 27
$ws.Range("D27").Value = "'1.669"
$ws.Range("E27").Value = "  -4.06%  "

# Row 28
$ws.Range("D28").Value = "'17.16"
$ws.Range("E28").Value = "  -5.21%  "

# Row 29
$ws.Range("D29").Value = "'109.31"
$ws.Range("E29").Value = "  -4.20%  "

# Row 30
$ws.Range("D30").Value = "'4.336"
$ws.Range("E30").Value = "  -8.03%  "

# Row 31
$ws.Range("D31").Value = "'4.237"
$ws.Range("E31").Value = "  -8.27%  "

# Row 32
$ws.Range("E32").Value = "  -3.73%  "

# Row 33
$ws.Range("D33").Value = "'0.04909"
$ws.Range("E33").Value = "  -2.16%  "

# Row 34
$ws.Range("D34").Value = "'0.7336"
$ws.Range("E34").Value = "  -9.78%  "

# Row 35
$ws.Range("D35").Value = "'1.139"
$ws.Range("E35").Value = "  -3.00%  "

# Row 36
$ws.Range("D36").Value = "'2.873"
$ws.Range("E36").Value = "  -3.05%  "

# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "'3.142"
$ws.Range("E37").Value = "  -1.99%  "

# Row 38
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'2.386"
$ws.Range("E38").Value = "  -9.99%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = "  -5.31%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5156"
$ws.Range("E40").Value = "  -15.43%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.9656"
$ws.Range("E41").Value = "  -10.01%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.232"
$ws.Range("E42").Value = "  -5.63%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'111.11"
$ws.Range("E43").Value = "  -3.43%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'8.057"
$ws.Range("E44").Value = "  -9.22%  "

# Row 45
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9999"
$ws.Range("E45").Value = "  -0.21%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4561"
$ws.Range("E46").Value = "  -11.79%  "

# Row 47
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.1371"
$ws.Range("E47").Value = "  -8.46%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'36.61"
$ws.Range("E48").Value = "  -2.53%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.237"
$ws.Range("E49").Value = "  -7.69%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.503"
$ws.Range("E50").Value = "  -8.39%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05821"
$ws.Range("E51").Value = "  -4.02%  "
